# Bug 50447: rename default colour scheme "Office" -> "New Office"
#
# The underlying change touches ppt/theme/theme1.xml: it swaps the
# accent1/accent5 theme colours and switches the Latin theme fonts from the
# Japanese "Yu Gothic" family to "Arial". (The source commit also relabels
# the <a:theme>/<a:clrScheme>/<a:fontScheme>/<a:fmtScheme> "name" attributes
# and edits the per-script Jpan/Hans <a:font> substitutions, but the
# PowerPoint object model has no writable property for those — Theme.Name,
# ThemeColorScheme.Name and ThemeFontScheme.Name are get-only/cosmetic here,
# and MajorFont/MinorFont only expose the 3 generic slots Latin/EastAsian/
# ComplexScript, not the extra per-script font table — so they are left as
# the host stores them.)
#
# We reach the theme through the slide master (every slide shares the one
# theme part in this deck).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$master = $s.Master
$theme = $master.Theme

# --- Colour scheme: swap accent1 <-> accent5 -------------------------------
# ThemeColorScheme.Item(n).RGB uses the VBA RGB() packing (R + G*256 + B*65536),
# i.e. the low byte is Red -- the reverse of the "RRGGBB" hex string stored in
# the OOXML <a:srgbClr val="RRGGBB"/>. Slots are ordered
# dk1,lt1,dk2,lt2,accent1..accent6,hlink,folHlink -> accent1=5, accent5=9.
$tcs = $theme.ThemeColorScheme
$tcs.Item(5).RGB = 0xD59B5B   # accent1: 4472C4 -> 5B9BD5
$tcs.Item(9).RGB = 0xC47244   # accent5: 5B9BD5 -> 4472C4

# --- Font scheme: Yu Gothic -> Arial for the Latin faces -------------------
$majorFont = $theme.ThemeFontScheme.MajorFont
$minorFont = $theme.ThemeFontScheme.MinorFont
$majorFont.Latin = "Arial"
$minorFont.Latin = "Arial"
